$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row from hunk 0
$ws.Range("H62").Value = 1978.5
$ws.Range("I62").Value = 2530.8333
$ws.Range("K62").Value = 2530.8333
$ws.Range("M62").Value = -1906.8333

# row from hunk 1
$ws.Range("H64").Value = 3303.1428
$ws.Range("I64").Value = 3285.7144
$ws.Range("J64").Value = 3311.8572
$ws.Range("K64").Value = 3285.7144
$ws.Range("L64").Value = 3311.8572
$ws.Range("M64").Value = -3037.7144
$ws.Range("N64").Value = -3807.8572

# row from hunk 2
$ws.Range("H65").Value = 1978.5
$ws.Range("I65").Value = 2530.8333
$ws.Range("K65").Value = 12654.1665
$ws.Range("M65").Value = -9534.166499999999

# row from hunk 3
$ws.Range("H67").Value = 3303.1428
$ws.Range("I67").Value = 3285.7144
$ws.Range("J67").Value = 3311.8572
$ws.Range("K67").Value = 3285.7144
$ws.Range("L67").Value = 3311.8572
$ws.Range("M67").Value = -2427.7144
$ws.Range("N67").Value = -5027.8572

# row from hunk 4
$ws.Range("H69").Value = 4161.364
$ws.Range("I69").Value = 4760
$ws.Range("J69").Value = 3662.5
$ws.Range("K69").Value = 14280
$ws.Range("L69").Value = 10987.5
$ws.Range("M69").Value = -13406
$ws.Range("N69").Value = -12735.5

# row from hunk 5
$ws.Range("H72").Value = 4161.364
$ws.Range("I72").Value = 4760
$ws.Range("J72").Value = 3662.5
$ws.Range("K72").Value = 42840
$ws.Range("L72").Value = 32962.5
$ws.Range("M72").Value = -38472
$ws.Range("N72").Value = -41698.5

# row from hunk 6
$ws.Range("H74").Value = 3877.2222
$ws.Range("I74").Value = 3997.3333
$ws.Range("J74").Value = 3853.2
$ws.Range("K74").Value = 3997.3333
$ws.Range("L74").Value = 3853.2
$ws.Range("M74").Value = -3061.3333
$ws.Range("N74").Value = -5725.2

# row from hunk 7
$ws.Range("H77").Value = 3877.2222
$ws.Range("I77").Value = 3997.3333
$ws.Range("J77").Value = 3853.2
$ws.Range("K77").Value = 19986.6665
$ws.Range("L77").Value = 19266
$ws.Range("M77").Value = -15306.6665
$ws.Range("N77").Value = -28626

# row from hunk 8
$ws.Range("H132").Value = 3455.6287
$ws.Range("I132").Value = 3636.0344
$ws.Range("J132").Value = 2583.6667
$ws.Range("K132").Value = 10908.1032
$ws.Range("L132").Value = 7751.000100000001
$ws.Range("M132").Value = -8378.1032
$ws.Range("N132").Value = -12811.0001

# row from hunk 9
$ws.Range("H138").Value = 2327.0715
$ws.Range("I138").Value = 1542.1111
$ws.Range("J138").Value = 3740
$ws.Range("K138").Value = 4626.3333
$ws.Range("L138").Value = 11220
$ws.Range("M138").Value = 513.6666999999998
$ws.Range("N138").Value = -21500

$ws = $wb.Worksheets.Item("ARM")
# row from hunk 10
$ws.Range("H32").Value = 7855.5425
$ws.Range("I32").Value = 8791.9
$ws.Range("J32").Value = 2653.5557
$ws.Range("K32").Value = 8791.9
$ws.Range("L32").Value = 2653.5557
$ws.Range("M32").Value = -8504.9
$ws.Range("N32").Value = -3227.5557

# row from hunk 11
$ws.Range("H45").Value = 2003.0454
$ws.Range("I45").Value = 1931.5
$ws.Range("J45").Value = 2325
$ws.Range("K45").Value = 1931.5
$ws.Range("L45").Value = 2325
$ws.Range("M45").Value = -1554.5
$ws.Range("N45").Value = -3079

# row from hunk 12
$ws.Range("H61").Value = 17243456
$ws.Range("I61").Value = 18520564
$ws.Range("K61").Value = 18520564
$ws.Range("M61").Value = -18520352

# row from hunk 13
$ws.Range("H136").Value = 17243456
$ws.Range("I136").Value = 18520564
$ws.Range("K136").Value = 55561692
$ws.Range("M136").Value = -55559142

$ws = $wb.Worksheets.Item("BSM")
# row from hunk 14
$ws.Range("H107").Value = 1872.5714
$ws.Range("I107").Value = 2112.5
$ws.Range("J107").Value = 433
$ws.Range("K107").Value = 2112.5
$ws.Range("L107").Value = 433
$ws.Range("M107").Value = -192.5
$ws.Range("N107").Value = -4273

$ws = $wb.Worksheets.Item("CRP")
# row from hunk 15
$ws.Range("H63").Value = 33271
$ws.Range("J63").Value = 33271
$ws.Range("L63").Value = 33271
$ws.Range("N63").Value = -34643

# row from hunk 16
$ws.Range("H66").Value = 33271
$ws.Range("J66").Value = 33271
$ws.Range("L66").Value = 99813
$ws.Range("N66").Value = -106677

# row from hunk 17
$ws.Range("H69").Value = 34700.6
$ws.Range("J69").Value = 50201
$ws.Range("L69").Value = 50201
$ws.Range("N69").Value = -51699

# row from hunk 18
$ws.Range("H72").Value = 34700.6
$ws.Range("J72").Value = 50201
$ws.Range("L72").Value = 150603
$ws.Range("N72").Value = -158091

# row from hunk 19
$ws.Range("H107").Value = 490.41177
$ws.Range("I107").Value = 565.0769
$ws.Range("J107").Value = 247.75
$ws.Range("K107").Value = 565.0769
$ws.Range("L107").Value = 247.75
$ws.Range("M107").Value = 1354.9231
$ws.Range("N107").Value = -4087.75

# row from hunk 20
$ws.Range("H140").Value = 30250
$ws.Range("J140").Value = 30250
$ws.Range("L140").Value = 30250
$ws.Range("N140").Value = -40610

$ws = $wb.Worksheets.Item("CUL")
# row from hunk 21
$ws.Range("H5").Value = 1382.1765
$ws.Range("I5").Value = 898.1429000000001
$ws.Range("K5").Value = 2694.4287
$ws.Range("M5").Value = -2582.4287

# row from hunk 22
$ws.Range("H95").Value = 11966.667
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 11966.667
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 35900.001
$ws.Range("N95").Value = -40018.001
$ws.Range("M95").Value = ""

# row from hunk 23
$ws.Range("H113").Value = 2091.3333
$ws.Range("I113").Value = 453.33334
$ws.Range("J113").Value = 3729.3333
$ws.Range("K113").Value = 1360.00002
$ws.Range("L113").Value = 11187.9999
$ws.Range("M113").Value = 809.9999800000001
$ws.Range("N113").Value = -15527.9999

# row from hunk 24
$ws.Range("H122").Value = 1391.8636
$ws.Range("I122").Value = 1137.2
$ws.Range("J122").Value = 1937.5714
$ws.Range("K122").Value = 10234.8
$ws.Range("L122").Value = 17438.1426
$ws.Range("M122").Value = -7784.800000000001
$ws.Range("N122").Value = -22338.1426

# row from hunk 25
$ws.Range("H135").Value = 1382.1765
$ws.Range("I135").Value = 898.1429000000001
$ws.Range("K135").Value = 8083.2861
$ws.Range("M135").Value = -5548.2861

$ws = $wb.Worksheets.Item("GSM")
# row from hunk 26
$ws.Range("H5").Value = 5004
$ws.Range("I5").Value = 5004
$ws.Range("K5").Value = 5004
$ws.Range("M5").Value = -4892

# row from hunk 27
$ws.Range("H7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = ""

# row from hunk 28
$ws.Range("H8").Value = 0
$ws.Range("J8").Value = 0
$ws.Range("L8").Value = 0
$ws.Range("N8").Value = ""

# row from hunk 29
$ws.Range("H70").Value = 11205.629
$ws.Range("I70").Value = 13581.272
$ws.Range("K70").Value = 13581.272
$ws.Range("M70").Value = -13311.272

# row from hunk 30
$ws.Range("H73").Value = 11205.629
$ws.Range("I73").Value = 13581.272
$ws.Range("K73").Value = 13581.272
$ws.Range("M73").Value = -12645.272

# row from hunk 31
$ws.Range("H132").Value = 5996.696
$ws.Range("I132").Value = 5184.4614
$ws.Range("J132").Value = 7052.6
$ws.Range("K132").Value = 15553.3842
$ws.Range("L132").Value = 21157.8
$ws.Range("M132").Value = -13023.3842
$ws.Range("N132").Value = -26217.8

# row from hunk 32
$ws.Range("H138").Value = 57839.6
$ws.Range("J138").Value = 57839.6
$ws.Range("L138").Value = 57839.6
$ws.Range("N138").Value = -68119.60000000001

$ws = $wb.Worksheets.Item("LTW")
# row from hunk 33
$ws.Range("H40").Value = 4463.52
$ws.Range("I40").Value = 3979.4
$ws.Range("J40").Value = 6400
$ws.Range("K40").Value = 3979.4
$ws.Range("L40").Value = 6400
$ws.Range("M40").Value = -3843.4
$ws.Range("N40").Value = -6672

# row from hunk 34
$ws.Range("H132").Value = 17252556
$ws.Range("I132").Value = 7840.3335
$ws.Range("J132").Value = 35729036
$ws.Range("K132").Value = 23521.0005
$ws.Range("L132").Value = 107187108
$ws.Range("M132").Value = -20991.0005
$ws.Range("N132").Value = -107192168

$ws = $wb.Worksheets.Item("WVR")
# row from hunk 35
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = ""

# row from hunk 36
$ws.Range("H124").Value = 14444
$ws.Range("J124").Value = 14444
$ws.Range("L124").Value = 14444
$ws.Range("N124").Value = -24264

# row from hunk 37
$ws.Range("H126").Value = 2172.8276
$ws.Range("I126").Value = 1696.174
$ws.Range("J126").Value = 4000
$ws.Range("K126").Value = 5088.522
$ws.Range("L126").Value = 12000
$ws.Range("M126").Value = -2618.522
$ws.Range("N126").Value = -16940
